# Apply the "PO Forecast" update:
#  1. Rename the "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
#  2. Rename the "Requested quantity" header on "Monthly Trend"   -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the existing "Requested quantity" headers -----------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" worksheet at the end ------------------------
$sheetCount = $wb.Worksheets.Count
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold + border + centered) from the Weekly sheet's header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-number-format cell style from the Weekly sheet's A2 onto column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A18").PasteSpecial(-4122)

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$forecastData = @(
    @(45557.99999999999, 40, -24.80072027385852, 107.2393504203942),
    @(45564.99999999999, 52, -11.95131903957503, 123.5229325261929),
    @(45571.99999999999, 63, -3.482643849463302, 128.7370574919426),
    @(45578.99999999999, 75, 5.200224685526235, 141.2034153421279),
    @(45585.99999999999, 86, 21.67055540620355, 156.2139575483153),
    @(45592.99999999999, 98, 29.43482817710758, 163.3736788413723),
    @(45599.99999999999, 109, 44.36492008965872, 178.1976326069847),
    @(45613.99999999999, 132, 63.70662158051771, 196.5599042616363),
    @(45620.99999999999, 144, 78.59946316224826, 211.2970028802509),
    @(45627.99999999999, 155, 90.16827546384532, 217.5485067457873),
    @(45634.99999999999, 167, 98.92311481389949, 235.8001436645033),
    @(45641.99999999999, 179, 115.8519686202036, 248.3844400672387),
    @(45648.99999999999, 190, 121.1230188286994, 260.0261678232119),
    @(45655.99999999999, 202, 139.2546680072611, 271.014930961211),
    @(45662.99999999999, 213, 147.428023391765, 283.4516533729741),
    @(45669.99999999999, 225, 162.6994048148129, 292.4995351334155),
    @(45676.99999999999, 236, 168.4558224573069, 299.6774898458108)
)

$row = 2
foreach ($rec in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

Write-Output "PO Forecast sheet added with $($row - 2) data rows"
